# Apply "cetak susulan dp2n15 dp1n9" edit: update cached MERGEFIELD result
# text for the NO / NAMA / SEPATU / TOPI fields (both copies of the form)
# and the UBN_1..UBN_9 / UH_1..UH_9 measurement fields.
#
# We resolve each field by its MERGEFIELD code name (not by searching for
# the old literal text) because several of the old numeric values repeat
# elsewhere in the document (e.g. "118" appears twice, "43" is both a
# UBN_9 value and the new SEPATU value), so a blind Find/Replace across
# the whole document would corrupt unrelated occurrences.

$d = $word.ActiveDocument

$map = @{
    "NO"     = "F47"
    "NAMA"   = "KADEK HARIAWAN"
    "SEPATU" = "43"
    "TOPI"   = "60"

    "UBN_1" = "50"
    "UBN_2" = "58"
    "UBN_3" = "21"
    "UBN_4" = "15"
    "UBN_5" = "120"
    "UBN_6" = "114"
    "UBN_7" = "122"
    "UBN_8" = "78"
    "UBN_9" = "48"

    "UH_1" = "50"
    "UH_2" = "59"
    "UH_3" = "20"
    "UH_4" = "13"
    "UH_5" = "30"
    "UH_6" = "29"
    "UH_7" = "30"
    "UH_8" = "77"
    "UH_9" = "48"
}

foreach ($f in $d.Fields) {
    $code = $f.Code.Text.Trim()
    if ($code.StartsWith("MERGEFIELD ")) {
        $name = $code.Substring(11).Trim()
        if ($map.ContainsKey($name)) {
            $r = $f.Result
            $d.Range($r.Start, $r.End).Text = $map[$name]
        }
    }
}
